# Generate Report for Handback
# Update the timestamp values recorded on the Overview/zh-cn/de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: Correspond Handback DateTime (G2) - Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-20 09:10:53"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-20 09:10:48"
$wsZhCn.Range("K2").Value = "2016-08-20 09:11:11"

# de-de sheet: Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-20 09:11:18"
